$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings (e.g. "1.002")
# are not auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '20.011.09'
$ws.Range("D3").Value = '1.422.93'
$ws.Range("E3").Value = '  -7.59%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '273.95'
$ws.Range("E6").Value = '  -5.62%  '
$ws.Range("D7").Value = '0.3754'
$ws.Range("E7").Value = '  -3.30%  '
$ws.Range("D8").Value = '0.3097'
$ws.Range("E8").Value = '  -2.89%  '
$ws.Range("D9").Value = '39.99'
$ws.Range("E9").Value = '  -7.72%  '
$ws.Range("D10").Value = '1.014'
$ws.Range("E10").Value = '  -4.45%  '
$ws.Range("D11").Value = '0.06599'
$ws.Range("E11").Value = '  -8.35%  '
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '5.403'
$ws.Range("E13").Value = '  -4.14%  '
$ws.Range("D14").Value = '17.13'
$ws.Range("E14").Value = '  -7.76%  '
$ws.Range("D15").Value = '6.171'
$ws.Range("E15").Value = '  -6.91%  '
$ws.Range("D16").Value = '1.425.95'
$ws.Range("E16").Value = '  -7.56%  '
$ws.Range("D17").Value = '0.00001011'
$ws.Range("E17").Value = '  -8.46%  '
$ws.Range("D18").Value = '0.05851'
$ws.Range("E18").Value = '  -11.03%  '
$ws.Range("D19").Value = '75.12'
$ws.Range("E19").Value = '  -9.99%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '5.661'
$ws.Range("E21").Value = '  -7.90%  '
$ws.Range("D22").Value = '14.49'
$ws.Range("E22").Value = '  -5.84%  '
$ws.Range("D23").Value = '11.00'
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").Value = '2.340'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").Value = '19.962.02'
$ws.Range("E25").Value = '  -8.28%  '
$ws.Range("D26").Value = '2.289'
$ws.Range("E26").Value = '  -4.33%  '
$ws.Range("D27").Value = '139.15'
$ws.Range("E27").Value = '  -4.55%  '
$ws.Range("D28").Value = '16.89'
$ws.Range("E28").Value = '  -8.11%  '
$ws.Range("D29").Value = '1.587.17'
$ws.Range("E29").Value = '  -7.61%  '
$ws.Range("D30").Value = '109.38'
$ws.Range("E30").Value = '  -7.02%  '
$ws.Range("D31").Value = '3.873'
$ws.Range("E31").Value = '  -20.07%  '
$ws.Range("D32").Value = '0.8920'
$ws.Range("E32").Value = '  -8.01%  '
$ws.Range("D33").Value = '5.434'
$ws.Range("E33").Value = '  -7.77%  '
$ws.Range("D34").Value = '0.07791'
$ws.Range("E34").Value = '  -5.20%  '
$ws.Range("D35").Value = '8.428'
$ws.Range("E35").Value = '  -5.98%  '
$ws.Range("D36").Value = '11.29'
$ws.Range("E36").Value = '  +5.96%  '
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '4.769'
$ws.Range("E38").Value = '  -7.09%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05704'
$ws.Range("E39").Value = '  -6.49%  '
$ws.Range("E40").Value = '  -6.28%  '
$ws.Range("D41").Value = '0.02024'
$ws.Range("E41").Value = '  -8.24%  '
$ws.Range("D42").Value = '1.094'
$ws.Range("E42").Value = '  -7.78%  '
$ws.Range("D43").Value = '1.269'
$ws.Range("E43").Value = '  -14.49%  '
$ws.Range("D44").Value = '0.5317'
$ws.Range("E44").Value = '  -7.53%  '
$ws.Range("D45").Value = '3.535'
$ws.Range("E45").Value = '  -5.66%  '
$ws.Range("D46").Value = '12.29'
$ws.Range("E46").Value = '  -5.68%  '
$ws.Range("D47").Value = '0.5137'
$ws.Range("E47").Value = '  -6.89%  '
$ws.Range("D48").Value = '1.790'
$ws.Range("E48").Value = '  -4.11%  '
$ws.Range("D49").Value = '109.91'
$ws.Range("D50").Value = '1.052'
$ws.Range("E50").Value = '  -7.74%  '
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.39%  '

# Restore default (unstyled) cell style now that text values are locked in,
# so the saved file does not carry a lingering explicit style index.
$ws.Range("D2:E51").Style = "Normal"
